# feat: add single and multi corrector
# Normalizes unit/range notation (x10E9/L -> 10^9/L, "~" -> "-", "fl" -> "fL"),
# trims trailing zeros in a few numeric result cells, and strips leading
# index numbers that had been accidentally prepended to two row labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to remain stored as text, preserving the exact
# string (avoids Excel auto-converting numeric-looking text into numbers).
function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# Row 2 - WBC
Set-TextValue "D2" "10^9/L"
Set-TextValue "F2" "3.5-9.5"

# Row 3 - RBC
Set-TextValue "D3" "10^12/L"
Set-TextValue "F3" "4.3-5.8"

# Row 4 - HGB
Set-TextValue "C4" "97.0"
Set-TextValue "F4" "130-175"

# Row 5 - PLT
Set-TextValue "D5" "10^9/L"
Set-TextValue "F5" "100-350"

# Row 6 - HCT
Set-TextValue "F6" "0.40-0.50"

# Row 7 - MCV
Set-TextValue "D7" "fL"
Set-TextValue "F7" "82-100"

# Row 8 - MCH
Set-TextValue "F8" "27-34"

# Row 9 - MCHC
Set-TextValue "F9" "316-354"

# Row 10 - NEUT
Set-TextValue "F10" "0.4-0.75"

# Row 11 - LYMPH
Set-TextValue "F11" "0.2-0.5"

# Row 12 - MONO
Set-TextValue "A12" "单核细胞百分率"
Set-TextValue "C12" "0.08"
Set-TextValue "F12" "0.030-0.100"

# Row 13 - EO
Set-TextValue "F13" "0.004-0.08"

# Row 14 - BASO
Set-TextValue "F14" "0-0.010"

# Row 15 - NEUT#
Set-TextValue "D15" "10^9/L"
Set-TextValue "F15" "1.8-6.3"

# Row 16 - LYMPH#
Set-TextValue "D16" "10^9/L"
Set-TextValue "F16" "1.1-3.2"

# Row 17 - MONO#
Set-TextValue "C17" "1.65"
Set-TextValue "D17" "10^9/L"
Set-TextValue "F17" "0.1-0.6"

# Row 18 - EO#
Set-TextValue "A18" "嗜酸性粒细胞绝对值"
Set-TextValue "C18" "5.61"
Set-TextValue "D18" "10^9/L"
Set-TextValue "F18" "0.02-0.52"

# Row 19 - BASO#
Set-TextValue "A19" "嗜碱性粒细胞绝对值"
Set-TextValue "C19" "0.16"
Set-TextValue "D19" "10^9/L"
Set-TextValue "F19" "0-0.06"
